# Commit: "Adicionando arquivos de PRA , e modificando outros"
#
# In this particular document, the reflective comment paragraph that was
# typed under "OUTRAS CONSIDERAÇÕES:" (praising "Professor Hugo") is
# removed. The paragraph mark / paragraph formatting (bold Arial rPr)
# stays in place — only its run content is deleted, leaving an empty
# paragraph, exactly like a user selecting the sentence and pressing
# Delete/Backspace in Word.

$d = $word.ActiveDocument

$commentText = "  Nesta situação onde o formador é extremamente profissional precisamos reconhecer e nos sentir privilegiados por ter o Professor Hugo desta formação, Além de suprema estou bem feliz, O professor Hugo é um otimo professor, e um ser humano admirável. Sucesso em qualquer lugar professor."

$rng = $d.Content
$found = $rng.Find.Execute($commentText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Clear the run text but keep the (now empty) paragraph and its pPr/rPr.
    $rng.Text = ""
} else {
    # Fallback: locate the paragraph via a unique marker substring and
    # clear its text (excluding the trailing paragraph mark), in case the
    # exact full-sentence match above ever misses (e.g. whitespace drift).
    $marker = "Professor Hugo desta formação"
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$marker*") {
            $clearRng = $d.Range($p.Range.Start, $p.Range.End - 1)
            $clearRng.Text = ""
            break
        }
    }
}
